$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 10 held a combined Akan/English sentence pair that is being
# split into two separate rows: the quote part, and the "fled" part.
# Insert one new row at row 11 so row 10 keeps the "quote" half and the new
# row 11 receives the "fled" half; everything that followed shifts down by one.
$ws.Rows("11:11").Insert()

$ws.Cells.Item(10, 1).Value = "Ye'see, ""Ananse eno no o? "" "
$ws.Cells.Item(10, 2).Value = "They said, ""Ananse, what about it? """

$ws.Cells.Item(11, 2).Value = " Ananse left that place; he has fled."
$ws.Cells.Item(11, 1).Value = "Ananse firii ho; wadwane."

$ws.Range("A11").Select()
